# Refer to resource name instead of jpg extension in spreadsheet.
# Strips the trailing ".jpg" from the "image on arrival" (column K) values
# and from a handful of Wikimedia/Wikipedia "image link" (column M) values
# on the "Sheet2" quest-stop table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Column K ("image on arrival"): strip ".jpg" row by row, top to bottom.
# Column M ("image link"): a few rows also end in ".jpg" (wikimedia/wikipedia
# links) and get the same treatment, right after that row's K cell.

$ws.Range("K3").Value = "stop_1"
$ws.Range("K4").Value = "stop_2"
$ws.Range("K5").Value = "stop_3"
$ws.Range("K6").Value = "stop_4"
$ws.Range("K7").Value = "stop_default"
$ws.Range("K8").Value = "stop_6"
$ws.Range("M8").Value = "https://commons.wikimedia.org/wiki/File:Waterhouse_Hylas_and_the_Nymphs_Manchester_Art_Gallery_1896.15"
$ws.Range("K9").Value = "stop_7"
$ws.Range("K10").Value = "stop_8"
$ws.Range("K11").Value = "stop_9"
$ws.Range("K12").Value = "stop_10"
$ws.Range("K13").Value = "stop_11"
$ws.Range("M13").Value = "https://commons.wikimedia.org/wiki/File:Clevelandart_1977.42"
$ws.Range("K14").Value = "stop_12"
$ws.Range("K15").Value = "stop_13"
$ws.Range("K16").Value = "stop_14"
$ws.Range("M16").Value = "https://commons.wikimedia.org/wiki/File:Peter_Paul_Rubens_007"
$ws.Range("K17").Value = "stop_15"
$ws.Range("K18").Value = "stop_16"
$ws.Range("K19").Value = "stop_18"
$ws.Range("M19").Value = "https://en.wikipedia.org/wiki/Argonautica#/media/File:MapoftheVoyageoftheArgonauts"

# Remaining rows all reuse the "stop_default" image.
$ws.Range("K20").Value = "stop_default"
$ws.Range("K21").Value = "stop_default"
$ws.Range("K22").Value = "stop_default"
$ws.Range("K23").Value = "stop_default"
$ws.Range("K24").Value = "stop_default"
$ws.Range("K25").Value = "stop_default"
$ws.Range("K26").Value = "stop_default"
$ws.Range("K27").Value = "stop_default"
$ws.Range("K28").Value = "stop_default"
$ws.Range("K29").Value = "stop_default"
$ws.Range("K30").Value = "stop_default"
$ws.Range("K31").Value = "stop_default"
$ws.Range("K32").Value = "stop_default"

# Leave the selection where the author ended up after editing the sheet.
$null = $ws.Range("A33").Select()
